# Updated DC Metro Data (Bus and Rail)
# - Rename headers to the new data-dictionary style names.
# - Re-key the YEAR-MONTH column from zero-padded "YYYY_MM" to "YYYY_M".
# - Drop the one-off "Coverage (miles)" figures that had been duplicated
#   across the 2012 rows.
# - Leave the cursor parked on the first blank row below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "YEAR-MONTH"
$ws.Range("B1").Value = "RIDERSHIP"
$ws.Range("C1").Value = "BUDGET "
$ws.Range("D1").Value = "Coverage (miles)"
$ws.Range("E1").Value = "POPULATION"

# --- Re-key column A ("YYYY_MM" -> "YYYY_M") ---------------------------
$lastRow = $ws.Range("A1").End(-4121).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value()
    $parts = $old.Split("_")
    $year = $parts[0]
    $month = [int]$parts[1]
    $cell.Value = "$year" + "_" + "$month"
}

# --- Clear the stray Coverage (miles) values on the 2012 rows ---------
$ws.Range("D38:D49").ClearContents()

# --- Park the selection where the author left it -----------------------
[void]$ws.Range("D74").Select()
